$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 272, pushing existing rows 272-290 down to 275-293.
$ws.Range("A272:A274").EntireRow.Insert()

# Fill in the 3 newly inserted rows with new weekly price data
# (columns A,B,C,E,F,G,N,Q,R repeat the same constant values used throughout this block)

# Row 272
$ws.Range("A272").Value = 4
$ws.Range("B272").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C272").Value = "Los Lagos"
$ws.Range("D272").Value = 44516
$ws.Range("E272").Value = 10
$ws.Range("F272").Value = 100112006
$ws.Range("G272").Value = "Repollo"
$ws.Range("H272").Value = "Copenhague"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 600
$ws.Range("K272").Value = 1400
$ws.Range("L272").Value = 1400
$ws.Range("M272").Value = 1400
$ws.Range("N272").Value = '$/unidad'
$ws.Range("O272").Value = "Región Metropolitana"
$ws.Range("P272").Value = 1400
$ws.Range("Q272").Value = 1
$ws.Range("R272").Value = "Hortaliza"

# Row 273
$ws.Range("A273").Value = 4
$ws.Range("B273").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C273").Value = "Los Lagos"
$ws.Range("D273").Value = 44516
$ws.Range("E273").Value = 10
$ws.Range("F273").Value = 100112006
$ws.Range("G273").Value = "Repollo"
$ws.Range("H273").Value = "Crespo record"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 500
$ws.Range("K273").Value = 1200
$ws.Range("L273").Value = 1200
$ws.Range("M273").Value = 1200
$ws.Range("N273").Value = '$/unidad'
$ws.Range("O273").Value = "Región Metropolitana"
$ws.Range("P273").Value = 1200
$ws.Range("Q273").Value = 1
$ws.Range("R273").Value = "Hortaliza"

# Row 274
$ws.Range("A274").Value = 4
$ws.Range("B274").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C274").Value = "Los Lagos"
$ws.Range("D274").Value = 44516
$ws.Range("E274").Value = 10
$ws.Range("F274").Value = 100112006
$ws.Range("G274").Value = "Repollo"
$ws.Range("H274").Value = "Crespo record"
$ws.Range("I274").Value = "Segunda"
$ws.Range("J274").Value = 500
$ws.Range("K274").Value = 1000
$ws.Range("L274").Value = 1000
$ws.Range("M274").Value = 1000
$ws.Range("N274").Value = '$/unidad'
$ws.Range("O274").Value = "Región Metropolitana"
$ws.Range("P274").Value = 1000
$ws.Range("Q274").Value = 1
$ws.Range("R274").Value = "Hortaliza"
